$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Johnny Hekker"
$ws.Range("B8").Value = "Group1"
$ws.Range("C8").Value = 45.73333333333333
$ws.Range("D8").Value = 64.33333333333333

$ws.Range("A9").Value = "Johnny Hekker"
$ws.Range("B9").Value = "Group2"
$ws.Range("C9").Value = 47
$ws.Range("D9").Value = 67.33333333333333

$ws.Range("A10").Value = "Johnny Hekker"
$ws.Range("B10").Value = "Difference"
$ws.Range("C10").Value = 1.266666666666673
$ws.Range("D10").Value = 3

$ws.Range("A11").Value = "Sam Martin"
$ws.Range("B11").Value = "Group1"
$ws.Range("C11").Value = 46.03333333333333
$ws.Range("D11").Value = 66.33333333333333

$ws.Range("A12").Value = "Sam Martin"
$ws.Range("B12").Value = "Group2"
$ws.Range("C12").Value = 46.73333333333333
$ws.Range("D12").Value = 67.33333333333333

$ws.Range("A13").Value = "Sam Martin"
$ws.Range("B13").Value = "Difference"
$ws.Range("C13").Value = 0.6999999999999957
$ws.Range("D13").Value = 1

$ws.Range("A2:D4").Copy()
$ws.Range("A8:D10").PasteSpecial(-4122)

$ws.Range("A5:D7").Copy()
$ws.Range("A11:D13").PasteSpecial(-4122)
